$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns hold values that *look* numeric ("215.52",
# "26.901.08", "0.531", ...) but are stored as plain text in the workbook
# (the site formats them itself, including thousands separators that are
# not valid numbers, e.g. "26.901.08"). Flip the whole data range to Text
# first so Excel's auto-type-detection on .Value assignment doesn't turn
# them into floating point numbers, then drop the number-format override
# again so the cells end up with no explicit style (matching how the rest
# of the sheet is built).
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# --- Price / Volume(1h) refresh for existing coin rows (rank unchanged) ---
$ws.Range("D2").Value = "26.901.08"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").Value = "1.668.69"
$ws.Range("E3").Value = "  +1.16%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "215.52"
$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("D6").Value = "0.531"
$ws.Range("E6").Value = "  +4.54%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  +1.23%  "

$ws.Range("E9").Value = "  +0.21%  "

$ws.Range("D10").Value = "20.28"
$ws.Range("E10").Value = "  +2.81%  "

$ws.Range("D11").Value = "0.0893"
$ws.Range("E11").Value = "  +3.27%  "

$ws.Range("D12").Value = "1.903.55"
$ws.Range("E12").Value = "  +1.13%  "

$ws.Range("D13").Value = "1.656.59"
$ws.Range("E13").Value = "  +0.36%  "

$ws.Range("E14").Value = "  +0.34%  "

$ws.Range("E15").Value = "  +1.82%  "

$ws.Range("D16").Value = "66.08"
$ws.Range("E16").Value = "  +1.50%  "

$ws.Range("D17").Value = "26.911.98"
$ws.Range("E17").Value = "  -0.39%  "

$ws.Range("D18").Value = "233.65"
$ws.Range("E18").Value = "  -1.84%  "

$ws.Range("D19").Value = "7.96"
$ws.Range("E19").Value = "  +1.63%  "

$ws.Range("E20").Value = "  +0.45%  "

$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("D22").Value = "4.40"
$ws.Range("E22").Value = "  -0.26%  "

$ws.Range("D23").Value = "9.12"
$ws.Range("E23").Value = "  -1.15%  "

$ws.Range("E24").Value = "  -2.76%  "

$ws.Range("D25").Value = "146.22"
$ws.Range("E25").Value = "  +0.52%  "

$ws.Range("D27").Value = "0.114"
$ws.Range("E27").Value = "  +0.65%  "

$ws.Range("D28").Value = "15.91"
$ws.Range("E28").Value = "  +0.68%  "

$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("E31").Value = "  +0.15%  "

$ws.Range("D32").Value = "3.35"
$ws.Range("E32").Value = "  +1.95%  "

$ws.Range("D33").Value = "1.454.95"
$ws.Range("E33").Value = "  -3.74%  "

$ws.Range("E34").Value = "  +2.37%  "

$ws.Range("D35").Value = "1.65"
$ws.Range("E35").Value = "  +3.90%  "

$ws.Range("E36").Value = "  -0.33%  "

$ws.Range("E37").Value = "  +1.05%  "

$ws.Range("D38").Value = "0.902"
$ws.Range("E38").Value = "  +2.02%  "

$ws.Range("E39").Value = "  +0.64%  "

$ws.Range("E40").Value = "  -3.50%  "

$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("E42").Value = "  +1.32%  "

$ws.Range("D43").Value = "66.13"
$ws.Range("E43").Value = "  +0.42%  "

$ws.Range("D44").Value = "0.975"
$ws.Range("E44").Value = "  +6.60%  "

$ws.Range("D45").Value = "1.812.58"
$ws.Range("E45").Value = "  +1.27%  "

$ws.Range("E46").Value = "  +1.41%  "

$ws.Range("D47").Value = "90.66"
$ws.Range("E47").Value = "  +1.45%  "

$ws.Range("E48").Value = "  +1.44%  "

# --- Bottom of the ranking reshuffled: BabyDogeCoin dropped off the list,
#     Algorand / Cronos move up a row, EnergySwap newly enters at #49 ---
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "0.102"
$ws.Range("E49").Value = "  +4.72%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.0506"
$ws.Range("E50").Value = "  -0.17%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.53"
$ws.Range("E51").Value = "  -0.01%  "

# Drop the temporary text-format override so the touched cells end up with
# no explicit cell style, same as every other data cell on the sheet.
$dataRange.ClearFormats()
